# Change from GDP to population as the default process emission driver
# (since this is continuous throughout)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sectors")

# Row 5 and rows 30-59 use "GDP" / "B2005USD" as the activity/units for
# process-emission-driven sectors ("NC" type). Replace the driver with
# population ("pop") measured in units of 1000.
$rows = @(5) + (30..59)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "pop"
    $ws.Cells.Item($r, 3).Value = 1000
}

# Update the view/selection to reflect where the edit was made.
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B30:C59").Select()
